# Apply the 4.3.1.1 workbook update:
#  - rename the English title string and add a 2021 ("M") column of data
#  - copy the 2020 ("L") column's number formatting so the new "M" column
#    (and the refreshed "L" column) look the same as the rest of the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the English sheet title (drop the period after "4.3.1.1").
#    Excel automatically re-packs the shared-string table, pruning the old
#    unused string and appending the new one, and re-numbers every other
#    <v> index that referenced a shared string after it.
$ws.Range("C1").Value2 = "4.3.1.1 Youth education by gender"

# 2. Add the new "2021" column (M) with the same look as the existing
#    "2020" column (L), which itself gets re-pasted so both columns share
#    identical formatting derived from the neighbouring "2019" column (K).
$ws.Range("K1:K12").Copy() | Out-Null
$ws.Range("L1:L12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("M1:M12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# 3. Header: year labels for the existing 2020 column and the new 2021 one.
$ws.Range("L3").Value2 = 2020
$ws.Range("M3").Value2 = 2021

# 4. Data rows: keep the 2020 ("L") values as-is and fill in the new 2021
#    ("M") figures.
$ws.Range("L4").Value2 = 10
$ws.Range("M4").Value2 = 10.8

$ws.Range("L5").Value2 = 6.4
$ws.Range("M5").Value2 = 5.2

$ws.Range("L6").Value2 = 13.5
$ws.Range("M6").Value2 = 16.2

$ws.Range("L7").Value2 = 24.3
$ws.Range("M7").Value2 = 24.2

$ws.Range("L8").Value2 = 27.8
$ws.Range("M8").Value2 = 27.6

$ws.Range("L9").Value2 = 20.9
$ws.Range("M9").Value2 = 20.9

$ws.Range("L10").Value2 = 26.7
$ws.Range("M10").Value2 = 28.5

$ws.Range("L11").Value2 = 28.4
$ws.Range("M11").Value2 = 29.7

$ws.Range("L12").Value2 = 25
$ws.Range("M12").Value2 = 27.5

# 5. Match the selection left behind in the saved file.
$ws.Range("O2").Select() | Out-Null
